$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 2739.2222
$ws.Range("I39").Value = 2373
$ws.Range("J39").Value = 3197
$ws.Range("K39").Value = 7119
$ws.Range("L39").Value = 9591
$ws.Range("M39").Value = -6823
$ws.Range("N39").Value = -10183
$ws.Range("H74").Value = 6428.2856
$ws.Range("I74").Value = 5833
$ws.Range("K74").Value = 5833
$ws.Range("M74").Value = -4897
$ws.Range("H76").Value = 4949.75
$ws.Range("J76").Value = 4599.6665
$ws.Range("L76").Value = 4599.6665
$ws.Range("N76").Value = -5229.6665
$ws.Range("H77").Value = 6428.2856
$ws.Range("I77").Value = 5833
$ws.Range("K77").Value = 29165
$ws.Range("M77").Value = -24485
$ws.Range("H79").Value = 4949.75
$ws.Range("J79").Value = 4599.6665
$ws.Range("L79").Value = 4599.6665
$ws.Range("N79").Value = -6783.6665
$ws.Range("H92").Value = 50262.9
$ws.Range("I92").Value = 55734.11
$ws.Range("J92").Value = 1022
$ws.Range("K92").Value = 55734.11
$ws.Range("L92").Value = 1022
$ws.Range("M92").Value = -54486.11
$ws.Range("N92").Value = -3518
$ws.Range("H132").Value = 2893.7192
$ws.Range("I132").Value = 2314.0557
$ws.Range("J132").Value = 13327.667
$ws.Range("K132").Value = 6942.1671
$ws.Range("L132").Value = 39983.001
$ws.Range("M132").Value = -4412.1671
$ws.Range("N132").Value = -45043.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 701.54285
$ws.Range("J2").Value = 968.5
$ws.Range("L2").Value = 968.5
$ws.Range("N2").Value = -1194.5
$ws.Range("H63").Value = 3313.5715
$ws.Range("I63").Value = 2739.2
$ws.Range("J63").Value = 4749.5
$ws.Range("K63").Value = 2739.2
$ws.Range("L63").Value = 4749.5
$ws.Range("M63").Value = -2053.2
$ws.Range("N63").Value = -6121.5
$ws.Range("H66").Value = 3313.5715
$ws.Range("I66").Value = 2739.2
$ws.Range("J66").Value = 4749.5
$ws.Range("K66").Value = 13696
$ws.Range("L66").Value = 23747.5
$ws.Range("M66").Value = -10264
$ws.Range("N66").Value = -30611.5
$ws.Range("H116").Value = 701.54285
$ws.Range("J116").Value = 968.5
$ws.Range("L116").Value = 968.5
$ws.Range("N116").Value = -5556.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 701.54285
$ws.Range("J3").Value = 968.5
$ws.Range("L3").Value = 968.5
$ws.Range("N3").Value = -1196.5
$ws.Range("H94").Value = 10980.308
$ws.Range("I94").Value = 12774.4
$ws.Range("K94").Value = 12774.4
$ws.Range("M94").Value = -12323.4
$ws.Range("H105").Value = 4467.391
$ws.Range("I105").Value = 4264.278
$ws.Range("K105").Value = 4264.278
$ws.Range("M105").Value = -2517.278

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 95328
$ws.Range("I133").Value = 98998
$ws.Range("J133").Value = 94716.336
$ws.Range("K133").Value = 98998
$ws.Range("L133").Value = 94716.336
$ws.Range("M133").Value = -96468
$ws.Range("N133").Value = -99776.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 4236.4
$ws.Range("I59").Value = 1727.6666
$ws.Range("K59").Value = 5182.9998
$ws.Range("M59").Value = -4642.9998
$ws.Range("H112").Value = 5812.0835
$ws.Range("J112").Value = 17500
$ws.Range("L112").Value = 52500
$ws.Range("N112").Value = -54716
$ws.Range("H123").Value = 2365.4
$ws.Range("I123").Value = 2365.4
$ws.Range("K123").Value = 7096.200000000001
$ws.Range("M123").Value = -4646.200000000001
$ws.Range("H131").Value = 2278664.5
$ws.Range("I131").Value = 2288.4546
$ws.Range("J131").Value = 3037456.5
$ws.Range("K131").Value = 6865.3638
$ws.Range("L131").Value = 9112369.5
$ws.Range("M131").Value = -1825.3638
$ws.Range("N131").Value = -9122449.5
$ws.Range("H134").Value = 2677.5
$ws.Range("I134").Value = 2050.2778
$ws.Range("K134").Value = 6150.8334
$ws.Range("M134").Value = -1080.8334
$ws.Range("H137").Value = 4090.1538
$ws.Range("I137").Value = 1141.75
$ws.Range("J137").Value = 5400.5557
$ws.Range("K137").Value = 3425.25
$ws.Range("L137").Value = 16201.6671
$ws.Range("M137").Value = 1674.75
$ws.Range("N137").Value = -26401.6671
$ws.Range("H138").Value = 1784.3334
$ws.Range("I138").Value = 1515
$ws.Range("J138").Value = 1999.8
$ws.Range("K138").Value = 4545
$ws.Range("L138").Value = 5999.4
$ws.Range("M138").Value = 595
$ws.Range("N138").Value = -16279.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3429
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 132315.75
$ws.Range("I132").Value = 202761.4
$ws.Range("J132").Value = 14906.333
$ws.Range("K132").Value = 608284.2
$ws.Range("L132").Value = 44718.999
$ws.Range("M132").Value = -605754.2
$ws.Range("N132").Value = -49778.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1075.7
$ws.Range("I16").Value = 722.4286
$ws.Range("J16").Value = 1900
$ws.Range("K16").Value = 722.4286
$ws.Range("L16").Value = 1900
$ws.Range("M16").Value = -552.4286
$ws.Range("N16").Value = -2240
$ws.Range("H22").Value = 38214
$ws.Range("J22").Value = 3416
$ws.Range("L22").Value = 3416
$ws.Range("N22").Value = -4006
$ws.Range("H27").Value = 38214
$ws.Range("J27").Value = 3416
$ws.Range("L27").Value = 3416
$ws.Range("N27").Value = -3630
$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30540
$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -31872
$ws.Range("H132").Value = 79065.44
$ws.Range("I132").Value = 112895.18
$ws.Range("J132").Value = 4640
$ws.Range("K132").Value = 338685.54
$ws.Range("L132").Value = 13920
$ws.Range("M132").Value = -336155.54
$ws.Range("N132").Value = -18980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 28441.2
$ws.Range("I41").Value = 33000
$ws.Range("J41").Value = 27301.5
$ws.Range("K41").Value = 33000
$ws.Range("L41").Value = 27301.5
$ws.Range("M41").Value = -32610
$ws.Range("N41").Value = -28081.5
$ws.Range("H132").Value = 117427
$ws.Range("I132").Value = 128056.17
$ws.Range("K132").Value = 384168.51
$ws.Range("M132").Value = -381638.51
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -99120
